# FP1-log-G47-xxxx-350-1201.xlsx
# Commit: "Generated timing waveforms for LogicUnit.vhd"
#
# Fills in the previously-blank activity-log row 32 (last 4 digits, date,
# start time, end time, description) and updates the sheet's scroll
# position / active selection to match where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New activity log entry on row 32 -------------------------------------
# last 4 digits
$ws.Cells.Item(32, 2).Value = 6977
# date (2020-04-03)
$ws.Cells.Item(32, 3).Value = 43924
# starttime (00:12)
$ws.Cells.Item(32, 4).Value = 0.0083333333333333332
# endtime (00:23)
$ws.Cells.Item(32, 5).Value = 0.015972222222222224
# description
$ws.Cells.Item(32, 7).Value = "Obtained timing waveforms and added them to Documentation as per instructions."

# --- Update the view: scroll to C13 and select F29 ------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 3
$ws.Range("F29").Select()
